# "About to refactor cresearcher to reduce repeated logic"
#
# Duplicate the "Mass simulation working" sheet (the latest timing-run
# sheet) into a new sheet "created c_did_get_hit" with a fresh timing run,
# plus a "% Improvement" column comparing the new run back to the old one.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Mass simulation working")

# Duplicate the sheet (keeps all values/formulas/styles/number formats)
# and place the copy immediately after the source sheet.
$src.Copy([System.Reflection.Missing]::Value, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "created c_did_get_hit"

# --- New timing numbers for the "created c_did_get_hit" run ---
$new.Range("B2").Value = 6.1669999999999998
$new.Range("C2").Value = 6.1609999999999996
$new.Range("D2").Value = 6.1269999999999998

$new.Range("B3").Value = 5.4649999999999999
$new.Range("C3").Value = 5.468
$new.Range("D3").Value = 5.4409999999999998

$new.Range("B4").Value = 0.69699999999999995
$new.Range("C4").Value = 0.68899999999999995
$new.Range("D4").Value = 0.68400000000000005

# New run date (one day later than the "Mass simulation working" run)
$new.Range("G2").Value = 41811

# --- Add "% Improvement" column comparing against the previous sheet ---
$new.Range("B1").Copy() | Out-Null
$new.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$new.Range("F1").Value = "% Improvement"

$new.Range("E2").Copy() | Out-Null
$new.Range("F2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$new.Range("F2").Formula = "=('Mass simulation working'!E2-'created c_did_get_hit'!E2)/'Mass simulation working'!E2"

$new.Application.CutCopyMode = $false

# Approximate the auto-fit column widths Excel produced for the new sheet.
$new.Columns.Item(1).ColumnWidth = 7
$new.Columns.Item(2).ColumnWidth = 5
$new.Columns.Item(3).ColumnWidth = 5
$new.Columns.Item(4).ColumnWidth = 5
$new.Columns.Item(5).ColumnWidth = 8
$new.Columns.Item(6).ColumnWidth = 14
$new.Columns.Item(7).ColumnWidth = 11
$new.Columns.Item(8).ColumnWidth = 44

$new.PageSetup.Orientation = 1   # xlPortrait

# The new sheet becomes the active tab; select its last-used cell.
$new.Activate() | Out-Null
$new.Range("F18").Select() | Out-Null

# The old sheet is no longer the active tab - its saved selection becomes
# the full data range instead of the single cell that was active before.
$src.Range("A1:H6").Select() | Out-Null
$new.Activate() | Out-Null
